$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before row 4, shifting existing rows 4-11 down to 5-12
$ws.Rows.Item(4).Insert()

# Populate the new row 4 with the climate change factor variable
$ws.Range("A4").Value = "General"
$ws.Range("B4").Value = "climate_change_factor_gnrl_hydropower_availability"
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = ""
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 0.5

$cols = @("J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS")
foreach ($col in $cols) {
    $ws.Range("$col`4").Value = 1
}
